# Update the "nota_view" column (J) on the sheet:
# every row whose current nota_view value is 5 becomes 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # column J = 10
    $val = $cell.Value2
    if ($val -ne $null -and [double]$val -eq 5) {
        $cell.Value = 4
    }
}
